# Automated Excel update — append 06-jul column to "Prix Spot" and a
# 2025-07-04 row to "Gaz" and "CO2".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": add column W (06-jul) with hourly prices
# ---------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Clone the header formatting of the previous date column (bold, centered,
# bordered) onto the new header cell before writing its text.
$wsSpot.Range("V1").Copy()
$wsSpot.Range("W1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsSpot.Range("W1").Value = "06-jul"

$spotValues = @(
    86.18000000000001,
    72.55,
    57.04,
    41.43,
    39.28,
    36.98,
    33.94,
    32.11,
    35.35,
    11.14,
    18.49,
    14.94,
    12.69,
    9.58,
    2.99,
    0.65,
    0.65,
    0.65,
    11.22,
    31.53,
    31.61,
    35.61,
    70.19,
    71.44
)

for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 23).Value = $spotValues[$i]
}

# ---------------------------------------------------------------------
# Sheet "Gaz": append row 20 (2025-07-04, 32.775)
# ---------------------------------------------------------------------
# The date column holds plain text (e.g. "2025-06-16"), not real dates, so
# force a text number format first to stop the "2025-07-04" literal from
# being auto-parsed into a date serial; then drop back to the (unstyled)
# Normal style so the cell matches its unstyled neighbours exactly.
$wsGaz = $wb.Worksheets.Item("Gaz")
$gazDateCell = $wsGaz.Range("A20")
$gazDateCell.NumberFormat = "@"
$gazDateCell.Value = "2025-07-04"
$gazDateCell.Style = "Normal"
$wsGaz.Range("B20").Value = 32.775

# ---------------------------------------------------------------------
# Sheet "CO2": append row 20 (2025-07-04, 70.92)
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$co2DateCell = $wsCo2.Range("A20")
$co2DateCell.NumberFormat = "@"
$co2DateCell.Value = "2025-07-04"
$co2DateCell.Style = "Normal"
$wsCo2.Range("B20").Value = 70.92
